$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, with value + formatting cloned from G1 (the
# existing rightmost header cell) so it matches the other header cells' style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add data values for the new Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
